# Add a new "2022-Q3" worksheet (fund-holdings detail) positioned between
# "总计" and "2021-Q1", and record its summary row on the "总计" sheet.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)   # "总计"
$ws2021  = $wb.Worksheets.Item(2)   # "2021-Q1" (existing, left untouched)

function Set-TextValue($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# ---------------------------------------------------------------------
# 1) Update the "总计" summary sheet: push the existing "2021-Q1" row
#    down to row 3 (carrying the row-2 A-column index style along with
#    it), then write the new "2022-Q3" summary into row 2.
# ---------------------------------------------------------------------
$wsTotal.Range("B3").Value = $wsTotal.Range("B2").Value2
$wsTotal.Range("C3").Value = $wsTotal.Range("C2").Value2
$wsTotal.Range("D3").Value = $wsTotal.Range("D2").Value2
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)   # xlPasteFormats
$wsTotal.Range("A3").Value = 1

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 0.1

# ---------------------------------------------------------------------
# 2) Create the new "2022-Q3" sheet by copying "总计" (so it inherits the
#    same header/number style + page margins), positioned right after
#    "总计", then overwrite its content with the 2022-Q3 fund-holdings
#    detail.
# ---------------------------------------------------------------------
$wsTotal.Copy($null, $wsTotal)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# Extend the header/index styling (copied from "总计") across the wider
# H-column table and the extra data rows this sheet needs.
$wsQ3.Range("D1").Copy()
$wsQ3.Range("E1:H1").PasteSpecial(-4122)   # xlPasteFormats
$wsQ3.Range("A2").Copy()
$wsQ3.Range("A3:A4").PasteSpecial(-4122)   # xlPasteFormats

$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"

# Row 2
$wsQ3.Range("A2").Value = 0
Set-TextValue $wsQ3 "B2" "006165"
Set-TextValue $wsQ3 "C2" "建信中证1000指数增强A"
Set-TextValue $wsQ3 "D2" "3.87"
Set-TextValue $wsQ3 "E2" "84.02"
Set-TextValue $wsQ3 "F2" "1.60"
Set-TextValue $wsQ3 "G2" "0.0619"
$wsQ3.Range("H2").Value = 3

# Row 3
$wsQ3.Range("A3").Value = 1
Set-TextValue $wsQ3 "B3" "006166"
Set-TextValue $wsQ3 "C3" "建信中证1000指数增强C"
Set-TextValue $wsQ3 "D3" "1.89"
Set-TextValue $wsQ3 "E3" "84.02"
Set-TextValue $wsQ3 "F3" "1.60"
Set-TextValue $wsQ3 "G3" "0.0302"
$wsQ3.Range("H3").Value = 3

# Row 4
$wsQ3.Range("A4").Value = 2
Set-TextValue $wsQ3 "B4" "013442"
Set-TextValue $wsQ3 "C4" "建信中证1000指数增强E"
Set-TextValue $wsQ3 "D4" "0.18"
Set-TextValue $wsQ3 "E4" "84.02"
Set-TextValue $wsQ3 "F4" "1.60"
Set-TextValue $wsQ3 "G4" "0.0029"
$wsQ3.Range("H4").Value = 3

# ---------------------------------------------------------------------
# 3) Restore the original active sheet ("总计") so workbook view state
#    matches the pre-edit selection.
# ---------------------------------------------------------------------
$wsTotal.Activate()
